$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - order chosen to reproduce original shared-string insertion order
$ws.Range("A1").Value = "TestScript"
$ws.Range("E1").Value = "Address"
$ws.Range("D1").Value = "Name"
$ws.Range("F1").Value = "City"
$ws.Range("G1").Value = "ZipCode"
$ws.Range("H1").Value = "State"
$ws.Range("I1").Value = "CreditCardNumber"
$ws.Range("J1").Value = "Month"
$ws.Range("K1").Value = "Year"
$ws.Range("L1").Value = "NameOnCard"
$ws.Range("A2").Value = "BookAFlight"
$ws.Range("B1").Value = "SourceCity"
$ws.Range("C1").Value = "DestinationCity"

# Data row (row 2)
$ws.Range("B2").Value = "Paris"
$ws.Range("C2").Value = "London"
$ws.Range("D2").Value = "Sai Audithya"
$ws.Range("E2").Value = "E9 DMV Skyland"
$ws.Range("F2").Value = "Nellore"
$ws.Range("G2").Value = 524004
$ws.Range("H2").Value = "Andhra"
# Leading apostrophe forces these numeric-looking values to be stored as
# text with a quote-prefix (matches the "quotePrefix" cell style).
$ws.Range("I2").Value = "'1111222233334444"
$ws.Range("J2").Value = "'01"
$ws.Range("K2").Value = "'2022"
$ws.Range("L2").Value = "Sai Audithya S"

# Copy the existing header formatting (bold font + yellow fill) onto the
# newly-added header cells G1:L1 so they match A1:F1.
$ws.Range("A1:F1").Copy()
$ws.Range("G1:L1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column widths (ColumnWidth uses character units; Excel pads by 5/6 of a
# character internally, so the stored OOXML width differs slightly from the
# value assigned here - this mirrors normal Excel autofit behaviour).
$ws.Columns.Item(1).ColumnWidth = 10.833333333333334
$ws.Columns.Item(2).ColumnWidth = 9.666666666666666
$ws.Columns.Item(3).ColumnWidth = 14.0
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666
$ws.Columns.Item(5).ColumnWidth = 14.333333333333334
$ws.Columns.Item(6).ColumnWidth = 6.833333333333333
$ws.Columns.Item(7).ColumnWidth = 7.5
$ws.Columns.Item(8).ColumnWidth = 6.666666666666667
$ws.Columns.Item(9).ColumnWidth = 17.166666666666668
$ws.Columns.Item(10).ColumnWidth = 6.166666666666667
$ws.Columns.Item(11).ColumnWidth = 4.166666666666667
$ws.Columns.Item(12).ColumnWidth = 12.666666666666666

# Remove the stale explicit selection range left over from before the edit.
$null = $ws.Range("A1").Select()
